$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("P1").Value = "time hrs"

# Row 2 fix + new formula cell
$ws.Range("H2").Value = -0.0008
$ws.Range("P2").Formula = "=O2/3600"

# Row 3 new formula cell (shared formula group anchor)
$ws.Range("P3").Formula = "=O3/3600"

# Row 4 fill in values
$ws.Range("B4").Value = 1.13
$ws.Range("C4").Value = 0.075
$ws.Range("D4").Value = -0.088
$ws.Range("E4").Value = 74
$ws.Range("F4").Value = 4.05
$ws.Range("G4").Value = -0.21
$ws.Range("H4").Value = -0.0008
$ws.Range("I4").Value = 49
$ws.Range("J4").Value = 0.22
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = 0.13
$ws.Range("M4").Value = 77
$ws.Range("O4").Value = 15600
$ws.Range("P4").Formula = "=O4/3600"

# Row 5 fill in values
$ws.Range("F5").Value = 3.32
$ws.Range("G5").Value = -0.24
$ws.Range("H5").Value = 0.0007
$ws.Range("I5").Value = 52
$ws.Range("J5").Value = 0.48
$ws.Range("K5").Value = 11
$ws.Range("L5").Value = 0.13
$ws.Range("M5").Value = 86

# Row 6 fill in values
$ws.Range("B6").Value = 1.5
$ws.Range("C6").Value = 0.1
$ws.Range("D6").Value = -0.024
$ws.Range("E6").Value = 88
$ws.Range("F6").Value = 3.32
$ws.Range("G6").Value = -0.24
$ws.Range("H6").Value = 0.0007
$ws.Range("I6").Value = 58
$ws.Range("J6").Value = 0.08
$ws.Range("K6").Value = 11
$ws.Range("L6").Value = 0.04
$ws.Range("M6").Value = 75
$ws.Range("N6").Value = 1500000
$ws.Range("O6").Value = 81300
$ws.Range("P6").Formula = "=O6/3600"

# Selection matches the diff (active cell N6)
[void]$ws.Range("N6").Select()

# Page setup orientation portrait
$ws.PageSetup.Orientation = 1
